# Actualizando archivo de tarifas -- se agregó el CANAL de SANLUIS
#
# Adds four new rate rows (67-70) to the PLAZAS_TARIFAS sheet for the
# SANLUIS channel (mirroring the HOR_INI/HOR_FIN/TARIFA_SPOT/TARIFA
# MENCION BLOQUEABLE layout used by every other plaza block), then
# updates the on-screen selection/scroll position on all three sheets
# to match where the author ended up after editing.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PLAZAS_TARIFAS")
$ws2 = $wb.Worksheets.Item("PLAZAS_CANALES")
$ws3 = $wb.Worksheets.Item("TARIFAS_NACIONALES")

# Copy the formatting (borders/style) of an existing 4-row plaza block
# down into the new rows so the added rows look identical to the rest
# of the table (style index 1 in the original file).
$ws1.Range("A2:E5").Copy()
$ws1.Range("A67:E70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 67 -- QUERETARO 06-14 tier, reused to seed the new block
$ws1.Range("A67").Value = "QUERETARO"
$ws1.Range("B67").Value = 6
$ws1.Range("C67").Value = 14
$ws1.Range("D67").Value = 1377.7
$ws1.Range("E67").Value = 7808.5

# Row 68 -- SANLUIS 14-18 tier
$ws1.Range("A68").Value = "SANLUIS"
$ws1.Range("B68").Value = 14
$ws1.Range("C68").Value = 18
$ws1.Range("D68").Value = 2559.9
$ws1.Range("E68").Value = 7808.5

# Row 69 -- SANLUIS 18-24 tier
$ws1.Range("A69").Value = "SANLUIS"
$ws1.Range("B69").Value = 18
$ws1.Range("C69").Value = 24
$ws1.Range("D69").Value = 3937.6
$ws1.Range("E69").Value = 7808.5

# Row 70 -- SANLUIS 00-06 tier
$ws1.Range("A70").Value = "SANLUIS"
$ws1.Range("B70").Value = 0
$ws1.Range("C70").Value = 6
$ws1.Range("D70").Value = 3937.6
$ws1.Range("E70").Value = 7808.5

# Leave the cursor/scroll position where the author last left it on
# each sheet.
[void]$ws1.Range("B78").Select()
[void]$ws2.Range("B9").Select()
[void]$ws3.Range("F16").Select()
